$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 830470
$ws.Range("C2").Value = "Norfolk Island"

# Row 3 (only amount changes)
$ws.Range("B3").Value = 9223

# Row 4
$ws.Range("B4").Value = 80749
$ws.Range("C4").Value = "Saint Lucia"

# Row 5
$ws.Range("B5").Value = 72605
$ws.Range("C5").Value = "Christmas Island"

# Row 6
$ws.Range("B6").Value = 734
$ws.Range("C6").Value = "Lao People's Democratic Republic"

# Row 7
$ws.Range("B7").Value = 1717
$ws.Range("C7").Value = "Portugal"

# Row 8
$ws.Range("B8").Value = 1160
$ws.Range("C8").Value = "Maldives"

# Row 9
$ws.Range("B9").Value = 3061
$ws.Range("C9").Value = "Liberia"

# Row 10
$ws.Range("B10").Value = 19
$ws.Range("C10").Value = "Czechia"

# Row 11
$ws.Range("B11").Value = 190
$ws.Range("C11").Value = "Nicaragua"

# Row 12
$ws.Range("B12").Value = 54
$ws.Range("C12").Value = "Micronesia, Federated States of"

# Row 13
$ws.Range("B13").Value = 1
$ws.Range("C13").Value = "Maldives"

# Row 14
$ws.Range("B14").Value = 0
$ws.Range("C14").Value = "Cuba"

# Row 15
$ws.Range("B15").Value = 0
$ws.Range("C15").Value = "Guernsey"

# Row 16
$ws.Range("B16").Value = 0
$ws.Range("C16").Value = "Svalbard and Jan Mayen"

# Row 17
$ws.Range("B17").Value = 0
$ws.Range("C17").Value = "Andorra"

# Row 18
$ws.Range("B18").Value = 17
$ws.Range("C18").Value = "Palau"
